# "back to original models"
# Revert the cheese-cave device model references / init-data JSON back to
# their original values (device model version ;1 instead of ;2, and the
# simplified alert-only JSON payloads instead of the humidity/temperature
# experiment values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# E5 reverts to the same "inUse: false" init data already used in E4
$ws.Range("E5").Value = '{"desiredTemperature": 50, "desiredHumidity": 85, "temperature": 50, "humidity": 50, "inUse": false, "temperatureAlert": false, "humidityAlert": false, "fanAlert": false}'

# Device model id goes back from ;2 to ;1 for all three cave devices
$ws.Range("A6").Value = "dtmi:com:contoso:digital_factory:cheese_factory:cheese_cave_device;1"
$ws.Range("A7").Value = "dtmi:com:contoso:digital_factory:cheese_factory:cheese_cave_device;1"
$ws.Range("A8").Value = "dtmi:com:contoso:digital_factory:cheese_factory:cheese_cave_device;1"

# Init data for the three cave devices reverts to the simple alert-only JSON
$deviceInitData = '{"desiredTemperature": 50, "desiredHumidity": 85, "temperatureAlert": true, "humidityAlert": true, "fanAlert": true}'
$ws.Range("E6").Value = $deviceInitData
$ws.Range("E7").Value = $deviceInitData
$ws.Range("E8").Value = $deviceInitData

# Restore selection to A6 as in the saved workbook
$ws.Activate()
$ws.Range("A6").Select()
